$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("H2").Value = 0.755982
$ws.Range("I2").Value = 0.3292948740286441
$ws.Range("J2").Value = 0.3292948740286441
$ws.Range("M2").Value = 1.660421
$ws.Range("N2").Value = 4.981262999999999
$ws.Range("O2").Value = 0.03714789785507311
$ws.Range("P2").Value = 0.03714789785507311
$ws.Range("Q2").Value = 0.4184161294739999
$ws.Range("R2").Value = 3.765745165266
$ws.Range("S2").Value = 0.01223261234461524
$ws.Range("T2").Value = 0.01223261234461524

# Row 3
$ws.Range("H3").Value = 0.755982
$ws.Range("I3").Value = 0.3292948740286441
$ws.Range("J3").Value = 0.3292948740286441
$ws.Range("O3").Value = 0.5631392661118858
$ws.Range("P3").Value = 0.5631392661118859
$ws.Range("R3").Value = 57.086378804142
$ws.Range("S3").Value = 0.1854388736948965
$ws.Range("T3").Value = 0.1854388736948966

# Row 4
$ws.Range("H4").Value = 0.755982
$ws.Range("I4").Value = 0.3292948740286441
$ws.Range("J4").Value = 0.3292948740286441
$ws.Range("M4").Value = 17.866195
$ws.Range("N4").Value = 53.598585
$ws.Range("O4").Value = 0.399712836033041
$ws.Range("P4").Value = 0.399712836033041
$ws.Range("Q4").Value = 4.50217394283
$ws.Range("R4").Value = 40.51956548547
$ws.Range("S4").Value = 0.1316233879891323
$ws.Range("T4").Value = 0.1316233879891323

# Row 5
$ws.Range("G5").Value = 0.5132593333333333
$ws.Range("H5").Value = 1.539778
$ws.Range("I5").Value = 0.6707051259713558
$ws.Range("J5").Value = 0.6707051259713558
$ws.Range("M5").Value = 1.660421
$ws.Range("N5").Value = 4.981262999999999
$ws.Range("O5").Value = 0.03714789785507311
$ws.Range("P5").Value = 0.03714789785507311
$ws.Range("Q5").Value = 0.8522265755126665
$ws.Range("R5").Value = 7.670039179613998
$ws.Range("S5").Value = 0.02491528551045787
$ws.Range("T5").Value = 0.02491528551045787

# Row 6
$ws.Range("G6").Value = 0.5132593333333333
$ws.Range("H6").Value = 1.539778
$ws.Range("I6").Value = 0.6707051259713558
$ws.Range("J6").Value = 0.6707051259713558
$ws.Range("O6").Value = 0.5631392661118858
$ws.Range("P6").Value = 0.5631392661118859
$ws.Range("Q6").Value = 12.91923032004644
$ws.Range("R6").Value = 116.273072880418
$ws.Range("S6").Value = 0.3777003924169892
$ws.Range("T6").Value = 0.3777003924169893

# Row 7
$ws.Range("G7").Value = 0.5132593333333333
$ws.Range("H7").Value = 1.539778
$ws.Range("I7").Value = 0.6707051259713558
$ws.Range("J7").Value = 0.6707051259713558
$ws.Range("M7").Value = 17.866195
$ws.Range("N7").Value = 53.598585
$ws.Range("O7").Value = 0.399712836033041
$ws.Range("P7").Value = 0.399712836033041
$ws.Range("Q7").Value = 9.169991334903333
$ws.Range("R7").Value = 82.52992201412999
$ws.Range("S7").Value = 0.2680894480439087
$ws.Range("T7").Value = 0.2680894480439087
